$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "TTC": the 256-simulation performance figures change from a
# constant 1200/1200/1200/2400 per-row block to a flat 870 everywhere,
# and the H-column "% " formulas are re-pointed from the O (1024) column
# to the M (256) column, matching "Performance average for 256
# simulations workflow".
# ---------------------------------------------------------------------
$ttc = $wb.Worksheets.Item("TTC")

# M:P values for rows 2, 4-11 all become 870 (row 3 has no M:P values).
$ttc.Range("M2:P2").Value = 870
$ttc.Range("M4:P4").Value = 870
$ttc.Range("M5:P5").Value = 870
$ttc.Range("M6:P6").Value = 870
$ttc.Range("M7:P7").Value = 870
$ttc.Range("M8:P8").Value = 870
$ttc.Range("M9:P9").Value = 870
$ttc.Range("M10:P10").Value = 870
$ttc.Range("M11:P11").Value = 870

# H4 switches from referencing column O to column M.
$ttc.Range("H4").Formula = "=(M4/B4)*100"

# H5:H7 carry the shared formula "(O./B.)*100" -> "(M./B.)*100"; assign
# across the same group of cells so the engine keeps them as one shared
# formula (H5 master, H6/H7 following).
$ttc.Range("H5:H7").Formula = "=(M5/B5)*100"

# H11 becomes a present-but-empty cell (same numeric style as H8:H10),
# it holds no value/formula.
$ttc.Range("H11").NumberFormat = "0"

# ---------------------------------------------------------------------
# Sheets "Tw" and "Te": selection becomes the range A1:B7 (no longer a
# single active cell on B4).
# ---------------------------------------------------------------------
$tw = $wb.Worksheets.Item("Tw")
$tw.Range("A1:B7").Select()

$te = $wb.Worksheets.Item("Te")
$te.Range("A1:B7").Select()

# ---------------------------------------------------------------------
# View: TTC becomes the workbook's active/selected sheet (was "plots"),
# with H4 (not H2) selected. Activating TTC last also drops
# "tabSelected" from "plots", which previously had it.
# ---------------------------------------------------------------------
$ttc.Activate()
$ttc.Range("H4").Select()
